$d = $word.ActiveDocument

# --- "Version 3" -> "Version 4" --------------------------------------------
$verPara = $d.Paragraphs(3).Range
$verPara.Find.Execute("3", $false, $false, $false, $false, $false, $true, 1, $false, "4", 2)

# --- cached DATE field result: "1/10/22 11:19 AM" -> "9/8/23 1:19 PM" ------
$datePara = $d.Paragraphs(4).Range
$datePara.Find.Execute("1/10/22 11:19 AM", $false, $false, $false, $false, $false, $true, 1, $false, "9/8/23 1:19 PM", 2)

# --- "ISO 21434" -> "ISO/SAE 21434" (keep bold) -----------------------------
$isoRng = $d.Content
$isoRng.Find.Execute("ISO 21434", $false, $false, $false, $false, $false, $true, 1, $false, "ISO/SAE 21434", 2)

# --- "Creative Commons Attribution-Share Alike (CC4-SA)" -------------------
#     -> "Creative Commons Attribution-Share Alike (CC-SA-4.0)" (keep bold)
$ccRng = $d.Content
$ccRng.Find.Execute("Creative Commons Attribution-Share Alike (CC4-SA)", $false, $false, $false, $false, $false, $true, 1, $false, "Creative Commons Attribution-Share Alike (CC-SA-4.0)", 2)

# --- add a first-page header with the AVCDL phase ID ------------------------
$sec = $d.Sections(1)
$sec.PageSetup.DifferentFirstPageHeaderFooter = 1
$firstHeader = $sec.Headers(2)
$tab = [char]9
$headerText = "$tab$tab" + "AVCDL-Supplier-3.1"
$hRng = $firstHeader.Range
$hRng.Collapse(0)
$hRng.InsertAfter($headerText)
$firstHeader.Range.Paragraphs(1).Style = "Header"

Write-Output "done"
